# App test bug log - add new test results / bug rows and fix statuses
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- status markers / newly tested results on existing rows ---
$ws.Range("F2").Value = "已解决"
$ws.Range("D6").Value = "已解决"
$ws.Range("F6").Value = "已解決"
$ws.Range("B7").Value = "已解決"
$ws.Range("D10").Value = "已解决"
$ws.Range("F10").Value = "已解決"
$ws.Range("D12").Value = "已解決"
$ws.Range("B14").Value = "已解决"

# --- new bug entries on rows 15-19 ---
$ws.Range("C15").Value = "15.班课成员太多，成员ListView 加载会头像错乱"

$ws.Range("B16").Value = "测试中"
$ws.Range("C16").Value = "16.作业名太长会导致经验值明细界面文字重叠"

$ws.Range("C17").Value = "17.姓名改成maxLength=""5"""

$ws.Range("C18").Value = "18.班课详情头像尽量少占用空间，右边第三行的数据显示不开，第三行三个字段的间隙也减少些"

$ws.Range("B19").Value = "已解决"
$ws.Range("C19").Value = "19.使用默认头像的用户刷新用户界面会一直刷新不结束"

# --- update view state (scrolled position / active selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("B18").Select()
